# Refreshes the cryptos list Price (D) and Volume(1h) (E) columns for
# rows 2-51, matching the Sat Sep 23 03:27:30 UTC 2023 GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.633.00"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "
$ws.Cells.Item(3, 4).Value = "1.596.70"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "
$ws.Cells.Item(4, 5).Value = "  +0.19%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.25"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.39%  "
$ws.Cells.Item(6, 5).Value = "  +0.53%  "
$ws.Cells.Item(7, 5).Value = "  +0.17%  "
$ws.Cells.Item(8, 5).Value = "  -0.04%  "
$ws.Cells.Item(9, 5).Value = "  -0.30%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.44"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.27%  "
$ws.Cells.Item(11, 5).Value = "  +0.05%  "
$ws.Cells.Item(12, 4).Value = "1.820.40"
$ws.Cells.Item(12, 5).Value = "  +0.07%  "
$ws.Cells.Item(13, 4).Value = "1.625.56"
$ws.Cells.Item(13, 5).Value = "  +1.94%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.03"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.14%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.522"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.51%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.88"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.45%  "
$ws.Cells.Item(17, 4).Value = "26.621.73"
$ws.Cells.Item(17, 5).Value = "  +0.01%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0736"
$ws.Cells.Item(18, 5).Value = "  +0.70%  "
$ws.Cells.Item(19, 5).Value = "  +0.30%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "208.47"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.56%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.04"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +4.98%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.27"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.20%  "
$ws.Cells.Item(23, 5).Value = "  -0.71%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.89"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.02%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "145.28"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.75%  "
$ws.Cells.Item(26, 5).Value = "  +0.13%  "
$ws.Cells.Item(27, 5).Value = "  -0.20%  "
$ws.Cells.Item(28, 5).Value = "  +0.18%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.24"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.73%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0507"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.17%  "
$ws.Cells.Item(31, 5).Value = "  -0.05%  "
$ws.Cells.Item(32, 5).Value = "  -0.45%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.93"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.66%  "
$ws.Cells.Item(34, 4).Value = "1.275.75"
$ws.Cells.Item(34, 5).Value = "  -1.59%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.46"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.38%  "
$ws.Cells.Item(37, 5).Value = "  -0.33%  "
$ws.Cells.Item(38, 5).Value = "  -0.73%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.839"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.59%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "5.49"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.25%  "
$ws.Cells.Item(41, 5).Value = "  +16.33%  "
$ws.Cells.Item(42, 5).Value = "  +0.63%  "
$ws.Cells.Item(43, 5).Value = "  -1.10%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "64.06"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.79%  "
$ws.Cells.Item(45, 4).Value = "1.732.60"
$ws.Cells.Item(45, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "90.19"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.93%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.60"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.37%  "
$ws.Cells.Item(48, 5).Value = "  +3.58%  "
$ws.Cells.Item(49, 5).Value = "  +1.10%  "
$ws.Cells.Item(50, 5).Value = "  +0.21%  "
$ws.Cells.Item(51, 5).Value = "  -1.11%  "
